# Update orch_session_issue_id (column H) UUID values in the 'diagnostics' sheet
# to newly generated UUIDs, reflecting the SQLa type-safety regeneration.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diagnostics")

$updates = @(
    @{Cell="H2"; Value="38749fe7-b3cd-4f6e-97a2-c7f4e3f8fdf2"},
    @{Cell="H3"; Value="61d9a7d6-09b4-4fa4-97c7-220513b0962f"},
    @{Cell="H4"; Value="72ee7100-3a87-4ffa-aef0-2640c5df8cfa"},
    @{Cell="H5"; Value="ad3b3796-a803-481a-a257-b321bbf6d1bb"},
    @{Cell="H6"; Value="ae8b6413-3631-4443-9d38-f9e073a456ad"},
    @{Cell="H7"; Value="e75ba306-36f4-45af-9a08-0a5e544bc4a0"},
    @{Cell="H8"; Value="60da1c0b-5068-4ac1-a523-891c5609d783"},
    @{Cell="H9"; Value="3d5a282e-4e68-4ad2-a5ab-61a4d9f51f36"},
    @{Cell="H10"; Value="17e73e9c-ed54-4f14-a303-35942c0bf1a4"},
    @{Cell="H11"; Value="112188d1-8d0a-46d8-81d5-a95331127302"},
    @{Cell="H12"; Value="424f0088-d833-465a-9e04-b82b9cef5847"},
    @{Cell="H13"; Value="29552173-a34c-4876-b3d2-68bbe72b9f03"},
    @{Cell="H14"; Value="8632494b-90c5-4ad5-a9f5-610bc116ff9a"},
    @{Cell="H15"; Value="da3d4b96-64f1-4852-9050-3c33c85d73a2"},
    @{Cell="H16"; Value="06ef4839-7e61-4c95-90ef-0290f03fd23b"},
    @{Cell="H17"; Value="27aeace6-65a0-4a33-8240-3694f7fa1f28"},
    @{Cell="H18"; Value="1b504d3f-18d7-482a-82c3-ac0ce2d4cc17"},
    @{Cell="H19"; Value="1a6b4a7e-dd5d-46b7-885b-fa0ccfac5470"},
    @{Cell="H20"; Value="82940c83-c0e9-4636-9e7b-9842532da0c0"},
    @{Cell="H21"; Value="6d9ff95f-3179-420b-b1c6-c0c5104e60e6"},
    @{Cell="H22"; Value="1bbf963b-1188-4e4d-91af-6fd31bfae417"},
    @{Cell="H27"; Value="117892cf-b527-401f-a916-1cd63b761a3e"},
    @{Cell="H28"; Value="73947d22-66c7-44fa-8b23-35c64dc10a23"},
    @{Cell="H29"; Value="4358499b-2b79-4251-a3b4-6297f7d901fc"},
    @{Cell="H30"; Value="0fb4bf12-5a94-4aa9-9481-0790070c4a7a"},
    @{Cell="H31"; Value="8208bb57-ff1a-4b53-8d71-e7313f90d59b"},
    @{Cell="H32"; Value="6014344d-23e6-41c3-8931-d3ae02cd786e"},
    @{Cell="H33"; Value="58b9ecad-c9e9-4358-9fd2-3debb8224e16"},
    @{Cell="H34"; Value="7bd902c7-513c-4514-acaa-bca929857511"},
    @{Cell="H35"; Value="4ad23d64-7ad5-49a8-bc68-c02073093772"},
    @{Cell="H36"; Value="20ad9557-2e53-4b45-8f8d-732fce01b81d"},
    @{Cell="H37"; Value="eec2748c-967b-4cd7-bdce-c6249409416d"},
    @{Cell="H38"; Value="dc8a448d-9e59-44ce-943b-d54d902b9ec5"},
    @{Cell="H40"; Value="a61ab194-a83f-4bcb-98dd-ad06b667e2bb"},
    @{Cell="H41"; Value="5b214c54-2e50-4e01-ae22-f0bb3e58a741"},
    @{Cell="H42"; Value="1034729c-f84d-486e-94fd-952cbd2fef7e"},
    @{Cell="H43"; Value="827fbffc-d9e7-403f-84fb-e152692ffd62"},
    @{Cell="H44"; Value="fdf0de5c-bc22-4b9d-9288-b51064314c88"},
    @{Cell="H45"; Value="067a5e3f-d138-4735-a96a-b0795c02af1a"},
    @{Cell="H46"; Value="788f51c2-447a-47b0-bb3c-89dfb0a183d5"},
    @{Cell="H47"; Value="8a4d4437-dcc9-474e-a21c-1fc98bdb9c39"},
    @{Cell="H48"; Value="8aa636d7-4e67-43cd-878a-8c670062131f"},
    @{Cell="H49"; Value="cc24ee55-84c2-410b-ae9b-e57f47975a3c"},
    @{Cell="H50"; Value="95d5d741-44eb-4c35-ba83-fd3a779823f5"},
    @{Cell="H51"; Value="01c5a162-4c1a-48a8-9a33-e02f5e5bdb9e"},
    @{Cell="H52"; Value="9432e199-654b-475c-b37c-6102f4a903e6"},
    @{Cell="H53"; Value="42b4cc12-b23d-4035-aa70-cb8fd1e023e2"},
    @{Cell="H54"; Value="ba0f9c27-ecf4-488f-9fbf-c53bc67a93f8"},
    @{Cell="H55"; Value="22c1e39d-d1e6-4836-a068-559369f60ff7"},
    @{Cell="H56"; Value="4e530af0-aa18-4d54-9ae6-d6b0107a5143"},
    @{Cell="H57"; Value="e667977b-8559-49eb-964c-ea8147517d71"},
    @{Cell="H58"; Value="4f8737d7-911d-41a1-8349-33247f90fb2c"},
    @{Cell="H59"; Value="9497e0a1-a4f9-4bce-8c17-aab924c62f08"},
    @{Cell="H60"; Value="fca586ac-95ae-4c20-a445-4178a2ff124d"},
    @{Cell="H61"; Value="a5d51571-6884-4e53-8603-f99816b808c2"},
    @{Cell="H62"; Value="85945493-6d00-4806-8a4a-0e19c11d14c1"},
    @{Cell="H63"; Value="51eeb4c9-f89c-44aa-a78b-8a85e1c4e0d1"},
    @{Cell="H64"; Value="1d531b7a-a791-4c50-98e1-d3afaf8edbc9"},
    @{Cell="H65"; Value="39372525-d66b-4f0d-8224-453c7b6e4956"},
    @{Cell="H66"; Value="c009beff-bc88-49cb-94ea-1431548b46b8"}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
